$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Каналы")

# --- Capture the existing comment-column (K) values before we start moving things around ---
$k4Value = $ws.Range("K4").Value2   # "сделать fixture.channel.count() из DB"
$k5Value = $ws.Range("K5").Value2   # "лекции про строки: 39-44"
$k6Value = $ws.Range("K6").Value2   # "тесты со случайными данными: 45"
$k7Value = $ws.Range("K7").Value2   # "DDT: 45, "
$k8Value = $ws.Range("K8").Value2   # "46 настройка параметры запуска в командной строке"

# --- Shift the lower block (old K5..K8) down two rows to K7..K10 ---
$ws.Range("K10").Value2 = $k8Value
$ws.Range("K9").Value2  = $k7Value
$ws.Range("K8").Value2  = $k6Value
$ws.Range("K7").Value2  = $k5Value

# --- Insert the two new "magic" comments. Write K6 first so the shared-string
#     table gets the two new unique strings in the same order as the target
#     workbook (index 24 = K6's text, index 25 = K5's text). ---
$ws.Range("K6").Value2 = "Магия pytest_generate_tests: 57 с 3:20"
$ws.Range("K5").Value2 = "Магия с фикстурой pytest: 18 с 08:28"

# --- Move the former K4 content up into K3 (it keeps the bold style already on K3) ---
$ws.Range("K3").Value2 = $k4Value

# --- Clear out the now-unused K4 cell entirely (no leftover value/formatting) ---
$ws.Range("K4").Clear()

# --- Update the saved selection to K6, matching the edited sheet view ---
$ws.Range("K6").Select()
